$wb = $excel.ActiveWorkbook

# --- BasicParsing sheet: insert new row 17 with "someNulls" / "_null_" test data ---
$ws1 = $wb.Worksheets.Item("BasicParsing")

$ws1.Rows.Item(17).Insert()
$ws1.Rows.Item(17).RowHeight = 18
$ws1.Range("A17").Value = "someNulls"
$ws1.Range("B17").Value = 1
$ws1.Range("C17").Value = "_null_"
$ws1.Range("D17").Value = $true
$ws1.Range("E17").Value = "_null_"
$ws1.Range("F17").Value = "hola"
$ws1.Range("G17").Value = '"all in"'
$ws1.Range("H17").Value = "_null_"

# widen column H to fit the new content
$ws1.Columns.Item(8).ColumnWidth = 35.166666666666664

# --- .EmptyRowTesting sheet: mirror the same sample row down at row 6 ---
$ws10 = $wb.Worksheets.Item(".EmptyRowTesting")
$ws1.Range("A17:H17").Copy($ws10.Range("A6:H6"))
$ws10.Rows.Item(6).RowHeight = 18

$ws10.Range("E16").Select()

# --- Selection / active sheet bookkeeping ---
$ws1.Activate()
$ws1.Range("A17:XFD17").Select()
